# "agregue el sistema de plantillas"
# Add a new row (id_audiencia=3) to the "audiencias" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("audiencias")

# Stage the new row's values in a scratch area, forced to Text so digit-
# looking / date-looking strings ("3", "123", "2023-05-11") aren't
# auto-converted to numbers/dates when assigned.
$ws.Range("H1:H6").NumberFormat = "@"
$ws.Range("H1").Value = "3"
$ws.Range("H2").Value = "123"
$ws.Range("H3").Value = "hueheu"
$ws.Range("H4").Value = "2023-05-11"
$ws.Range("H5").Value = "das"
$ws.Range("H6").Value = "dsfa"

# Copy just the values (not the Text number format) into row 4 so the
# destination cells keep the workbook's default (General) style while
# still being stored as text.
$ws.Range("H1").Copy()
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("H2").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("H3").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("H4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("H5").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("H6").Copy()
$ws.Range("F4").PasteSpecial(-4163)

# Clean up the scratch area.
$ws.Range("H1:H6").Clear()
